$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Rename the header row: "_old" columns become "_FV2404", "_new" columns
#    become "_FV2410" (the "diff" column in between stays unchanged).
# ---------------------------------------------------------------------------
$oldSuffixHeaders = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

# Columns A-J -> "<name>_FV2404"
for ($i = 0; $i -lt $oldSuffixHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = "$($oldSuffixHeaders[$i])_FV2404"
}

# Column K -> "diff" (unchanged)
$ws.Cells.Item(1, 11).Value = "diff"

# Columns L-U -> "<name>_FV2410"
for ($i = 0; $i -lt $oldSuffixHeaders.Length; $i++) {
    $ws.Cells.Item(1, 12 + $i).Value = "$($oldSuffixHeaders[$i])_FV2410"
}

# ---------------------------------------------------------------------------
# 2. Turn the data range into an Excel Table ("Table1") spanning A1:U65.
# ---------------------------------------------------------------------------
$tableRange = $ws.Range("A1:U65")
$table = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$table.Name = "Table1"

# ---------------------------------------------------------------------------
# 3. Freeze the header row (split/freeze at row 2).
# ---------------------------------------------------------------------------
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
